$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45898
$ws.Range("B2").Value = 82.7
$ws.Range("C2").Value = 79.64
$ws.Range("D2").Value = 72.2
$ws.Range("E2").Value = 60
$ws.Range("F2").Value = 54.43
$ws.Range("G2").Value = 54.6
$ws.Range("H2").Value = 75.01000000000001
$ws.Range("I2").Value = 80.2
$ws.Range("J2").Value = 70.09
$ws.Range("K2").Value = 42.47
$ws.Range("L2").Value = 5
$ws.Range("M2").Value = 0.83
$ws.Range("N2").Value = 0.65
$ws.Range("O2").Value = 0.01
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 3.52
$ws.Range("U2").Value = 53.55
$ws.Range("V2").Value = 81.8
$ws.Range("W2").Value = 100.01
$ws.Range("X2").Value = 95.14
$ws.Range("Y2").Value = 86.17
$ws.Range("Z2").Value = 45.75
$ws.Range("AB2").Value = 90.78
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 90.90000000000001
$ws.Range("AE2").Value = "22h-24h"
$ws.Range("AF2").Value = 90.66
$ws.Range("AG2").Value = "9h-18h"
